$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking values
# (e.g. "228.41") are stored as strings, matching the source data which
# uses inline/shared strings for every Price cell.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "41.642.68"
$ws.Range("E2").Value = "  +4.93%  "

# Row 3
$ws.Range("D3").Value = "2.225.48"
$ws.Range("E3").Value = "  +3.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").Value = "228.41"
$ws.Range("E5").Value = "  +0.48%  "

# Row 6
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  -0.90%  "

# Row 7
$ws.Range("D7").Value = "61.04"
$ws.Range("E7").Value = "  -3.60%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").Value = "0.401"
$ws.Range("E9").Value = "  +2.26%  "

# Row 10
$ws.Range("D10").Value = "58.01"
$ws.Range("E10").Value = "  -0.41%  "

# Row 11
$ws.Range("D11").Value = "0.0875"
$ws.Range("E11").Value = "  +3.45%  "

# Row 12
$ws.Range("E12").Value = "  -0.03%  "

# Row 13
$ws.Range("D13").Value = "2.554.22"
$ws.Range("E13").Value = "  +2.97%  "

# Row 14
$ws.Range("D14").Value = "15.66"
$ws.Range("E14").Value = "  -1.58%  "

# Row 15
$ws.Range("D15").Value = "21.52"
$ws.Range("E15").Value = "  -1.53%  "

# Row 16
$ws.Range("D16").Value = "0.794"
$ws.Range("E16").Value = "  -1.22%  "

# Row 17
$ws.Range("D17").Value = "5.55"
$ws.Range("E17").Value = "  +1.22%  "

# Row 18
$ws.Range("D18").Value = "2.228.50"
$ws.Range("E18").Value = "  +3.17%  "

# Row 19
$ws.Range("D19").Value = "41.555.63"
$ws.Range("E19").Value = "  +4.83%  "

# Row 20
$ws.Range("D20").Value = "72.55"
$ws.Range("E20").Value = "  +1.21%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0885"
$ws.Range("E21").Value = "  +4.94%  "

# Row 22
$ws.Range("D22").Value = "6.02"
$ws.Range("E22").Value = "  -0.70%  "

# Row 23
$ws.Range("D23").Value = "247.01"
$ws.Range("E23").Value = "  +7.43%  "

# Row 24
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  -0.81%  "

# Row 26
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  -1.88%  "

# Row 27
$ws.Range("D27").Value = "9.45"
$ws.Range("E27").Value = "  -0.61%  "

# Row 28
$ws.Range("D28").Value = "167.79"
$ws.Range("E28").Value = "  -2.64%  "

# Row 29
$ws.Range("E29").Value = "  -0.24%  "

# Row 30
$ws.Range("D30").Value = "19.88"
$ws.Range("E30").Value = "  +0.24%  "

# Row 31
$ws.Range("E31").Value = "  -3.22%  "

# Row 32
$ws.Range("D32").Value = "2.62"
$ws.Range("E32").Value = "  -2.43%  "

# Row 33
$ws.Range("D33").Value = "0.122"
$ws.Range("E33").Value = "  -0.52%  "

# Row 34
$ws.Range("E34").Value = "  +7.77%  "

# Row 35
$ws.Range("D35").Value = "4.64"
$ws.Range("E35").Value = "  +1.16%  "

# Row 36
$ws.Range("D36").Value = "0.0621"
$ws.Range("E36").Value = "  +0.32%  "

# Row 37
$ws.Range("E37").Value = "  -5.41%  "

# Row 38
$ws.Range("D38").Value = "3.69"
$ws.Range("E38").Value = "  +2.71%  "

# Row 39
$ws.Range("D39").Value = "2.37"
$ws.Range("E39").Value = "  -1.05%  "

# Row 40
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.13%  "

# Row 41
$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").Value = "0.000235"
$ws.Range("E41").Value = "  +28.42%  "

# Row 42
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "4.85"
$ws.Range("E42").Value = "  -5.41%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0235"
$ws.Range("E43").Value = "  +4.25%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "8.77"
$ws.Range("E44").Value = "  +13.07%  "

# Row 45
$ws.Range("D45").Value = "0.0988"
$ws.Range("E45").Value = "  +7.51%  "

# Row 46
$ws.Range("D46").Value = "98.99"
$ws.Range("E46").Value = "  -3.68%  "

# Row 47
$ws.Range("D47").Value = "1.469.52"
$ws.Range("E47").Value = "  -2.97%  "

# Row 48
$ws.Range("E48").Value = "  -2.30%  "

# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "16.32"
$ws.Range("E49").Value = "  -7.14%  "

# Row 50
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.78"
$ws.Range("E50").Value = "  -1.07%  "

# Row 51
$ws.Range("D51").Value = "1.08"
$ws.Range("E51").Value = "  -1.12%  "

# Restore the default (General/no explicit number-format) style on column D
# so the saved styles table matches the original (only the shared-string
# table / inline text changes, no new number formats are introduced).
$ws.Range("D2:D51").Style = "Normal"
